# Correct NGSI-LD attribute names
# Insert a new "Cleaned_KKS_name" column before the existing "OpenAPI_3.0_Type"
# column (old column N), shifting it (and the following "X-NGSI-LD" column)
# one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at N; this pushes the old N (OpenAPI_3.0_Type) and
# O (X-NGSI-LD) columns one position to the right (O and P).
$ws.Range("N1").EntireColumn.Insert()

# Find the last used row in column A (data rows 2..44).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Header for the newly inserted column, matching the style of the other
# header cells (bold / bordered header style already present on N1 before
# the shift put it here; copy style from the now-adjacent header cell).
$ws.Range("N1").Value = "Cleaned_KKS_name"
$ws.Range("O1").Copy()
$ws.Range("N1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

for ($r = 2; $r -le $lastRow; $r++) {
    $name = [string]$ws.Cells.Item($r, 1).Value2
    if ($name -ne $null) {
        $cleaned = $name.Replace("/", ":").Replace(".", ":")
        $ws.Cells.Item($r, 14).Value = $cleaned
    }
}

$wb.Save()
